# Correction in SA algorithm: update the logged Fitness values (column C)
# for run_15 with the corrected figures from the re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C30").Value = 8045
$ws.Range("C31:C49").Value = 7736
$ws.Range("C50:C79").Value = 7345
$ws.Range("C80:C94").Value = 7312
$ws.Range("C95:C103").Value = 7295
$ws.Range("C110:C143").Value = 7295
